$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.131.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.98%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.317.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.25%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.52"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.36%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.311.78"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.23%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.66%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.41%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.571"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.03"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.845.61"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.02%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.53%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "595.72"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.67%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.116.28"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.01%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.322.93"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.35%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.73"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.36%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.89"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.76%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.42%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.89"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.97"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.29%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.78"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -7.58%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.60%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.04%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.65%  "

# Row 29
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.45"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.47%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.67"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.31%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.53%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "560.52"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.00%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.64"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.19%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.86"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.03%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.804.45"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.51%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.59%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.07"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.84%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.22%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.06"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.26%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.39%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.40"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.90%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.66%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.17%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.06"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -8.54%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.66%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.05%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.51"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.53%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.80"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.67%  "
